$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.442.52'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '2.332.06'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.578'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '2.329.84'
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('D15').Value = '2.747.34'
$ws.Range('E15').Value = '  +1.01%  '
$ws.Range('D16').Value = '60.383.67'
$ws.Range('E16').Value = '  +3.08%  '
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '2.339.45'
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '315.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').Value = '  +5.56%  '
$ws.Range('E29').Value = '  +9.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('E34').Value = '  +9.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.380'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '322.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.76%  '
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.86%  '
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0939'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.07%  '
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('D50').Value = '0.0₆0215'
$ws.Range('E50').Value = '  +15.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.50%  '
